$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 1.987640333333333
$ws.Cells.Item(2, 8).Value = 5.962921
$ws.Cells.Item(2, 9).Value = 0.0278174819837782
$ws.Cells.Item(2, 10).Value = 0.0278174819837782
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 218.7785543333333
$ws.Cells.Item(2, 14).Value = 656.3356630000001
$ws.Cells.Item(2, 15).Value = 0.7837094150017259
$ws.Cells.Item(2, 16).Value = 0.7837094150017259
$ws.Cells.Item(2, 17).Value = 434.8530786612914
$ws.Cells.Item(2, 18).Value = 3913.677707951623
$ws.Cells.Item(2, 19).Value = 0.02180082253232786
$ws.Cells.Item(2, 20).Value = 0.02180082253232786

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 1.987640333333333
$ws.Cells.Item(3, 8).Value = 5.962921
$ws.Cells.Item(3, 9).Value = 0.0278174819837782
$ws.Cells.Item(3, 10).Value = 0.0278174819837782
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 46.29469433333333
$ws.Cells.Item(3, 14).Value = 138.884083
$ws.Cells.Item(3, 15).Value = 0.1658370397602197
$ws.Cells.Item(3, 16).Value = 0.1658370397602197
$ws.Cells.Item(3, 17).Value = 92.01720167627144
$ws.Cells.Item(3, 18).Value = 828.154815086443
$ws.Cells.Item(3, 19).Value = 0.00461316886577302
$ws.Cells.Item(3, 20).Value = 0.00461316886577302

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 1.987640333333333
$ws.Cells.Item(4, 8).Value = 5.962921
$ws.Cells.Item(4, 9).Value = 0.0278174819837782
$ws.Cells.Item(4, 10).Value = 0.0278174819837782
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 6.978882
$ws.Cells.Item(4, 14).Value = 20.936646
$ws.Cells.Item(4, 15).Value = 0.02499977909741928
$ws.Cells.Item(4, 16).Value = 0.02499977909741927
$ws.Cells.Item(4, 17).Value = 13.871507344774
$ws.Cells.Item(4, 18).Value = 124.843566102966
$ws.Cells.Item(4, 19).Value = 0.0006954309046408956
$ws.Cells.Item(4, 20).Value = 0.0006954309046408955

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 1.987640333333333
$ws.Cells.Item(5, 8).Value = 5.962921
$ws.Cells.Item(5, 9).Value = 0.0278174819837782
$ws.Cells.Item(5, 10).Value = 0.0278174819837782
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 7.105616
$ws.Cells.Item(5, 14).Value = 21.316848
$ws.Cells.Item(5, 15).Value = 0.02545376614063513
$ws.Cells.Item(5, 16).Value = 0.02545376614063513
$ws.Cells.Item(5, 17).Value = 14.12340895477867
$ws.Cells.Item(5, 18).Value = 127.110680593008
$ws.Cells.Item(5, 19).Value = 0.0007080596810364214
$ws.Cells.Item(5, 20).Value = 0.0007080596810364213

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 56.81334033333334
$ws.Cells.Item(6, 8).Value = 170.440021
$ws.Cells.Item(6, 9).Value = 0.7951157181995667
$ws.Cells.Item(6, 10).Value = 0.7951157181995667
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 218.7785543333333
$ws.Cells.Item(6, 14).Value = 656.3356630000001
$ws.Cells.Item(6, 15).Value = 0.7837094150017259
$ws.Cells.Item(6, 16).Value = 0.7837094150017259
$ws.Cells.Item(6, 17).Value = 12429.54046497433
$ws.Cells.Item(6, 18).Value = 111865.8641847689
$ws.Cells.Item(6, 19).Value = 0.6231396743688595
$ws.Cells.Item(6, 20).Value = 0.6231396743688595

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 56.81334033333334
$ws.Cells.Item(7, 8).Value = 170.440021
$ws.Cells.Item(7, 9).Value = 0.7951157181995667
$ws.Cells.Item(7, 10).Value = 0.7951157181995667
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 46.29469433333333
$ws.Cells.Item(7, 14).Value = 138.884083
$ws.Cells.Item(7, 15).Value = 0.1658370397602197
$ws.Cells.Item(7, 16).Value = 0.1658370397602197
$ws.Cells.Item(7, 17).Value = 2630.156224787305
$ws.Cells.Item(7, 18).Value = 23671.40602308574
$ws.Cells.Item(7, 19).Value = 0.1318596369730372
$ws.Cells.Item(7, 20).Value = 0.1318596369730372

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 56.81334033333334
$ws.Cells.Item(8, 8).Value = 170.440021
$ws.Cells.Item(8, 9).Value = 0.7951157181995667
$ws.Cells.Item(8, 10).Value = 0.7951157181995667
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 6.978882
$ws.Cells.Item(8, 14).Value = 20.936646
$ws.Cells.Item(8, 15).Value = 0.02499977909741928
$ws.Cells.Item(8, 16).Value = 0.02499977909741927
$ws.Cells.Item(8, 17).Value = 396.493598212174
$ws.Cells.Item(8, 18).Value = 3568.442383909566
$ws.Cells.Item(8, 19).Value = 0.01987771731187505
$ws.Cells.Item(8, 20).Value = 0.01987771731187504

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 56.81334033333334
$ws.Cells.Item(9, 8).Value = 170.440021
$ws.Cells.Item(9, 9).Value = 0.7951157181995667
$ws.Cells.Item(9, 10).Value = 0.7951157181995667
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 7.105616
$ws.Cells.Item(9, 14).Value = 21.316848
$ws.Cells.Item(9, 15).Value = 0.02545376614063513
$ws.Cells.Item(9, 16).Value = 0.02545376614063513
$ws.Cells.Item(9, 17).Value = 403.6937800859787
$ws.Cells.Item(9, 18).Value = 3633.244020773808
$ws.Cells.Item(9, 19).Value = 0.02023868954579492
$ws.Cells.Item(9, 20).Value = 0.02023868954579491

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 0.7501196666666666
$ws.Cells.Item(10, 8).Value = 2.250359
$ws.Cells.Item(10, 9).Value = 0.01049809664416703
$ws.Cells.Item(10, 10).Value = 0.01049809664416703
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 218.7785543333333
$ws.Cells.Item(10, 14).Value = 656.3356630000001
$ws.Cells.Item(10, 15).Value = 0.7837094150017259
$ws.Cells.Item(10, 16).Value = 0.7837094150017259
$ws.Cells.Item(10, 17).Value = 164.1100962503352
$ws.Cells.Item(10, 18).Value = 1476.990866253017
$ws.Cells.Item(10, 19).Value = 0.008227457179631727
$ws.Cells.Item(10, 20).Value = 0.008227457179631728

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 0.7501196666666666
$ws.Cells.Item(11, 8).Value = 2.250359
$ws.Cells.Item(11, 9).Value = 0.01049809664416703
$ws.Cells.Item(11, 10).Value = 0.01049809664416703
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 46.29469433333333
$ws.Cells.Item(11, 14).Value = 138.884083
$ws.Cells.Item(11, 15).Value = 0.1658370397602197
$ws.Cells.Item(11, 16).Value = 0.1658370397602197
$ws.Cells.Item(11, 17).Value = 34.72656068175522
$ws.Cells.Item(11, 18).Value = 312.539046135797
$ws.Cells.Item(11, 19).Value = 0.001740973270585357
$ws.Cells.Item(11, 20).Value = 0.001740973270585357

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 0.7501196666666666
$ws.Cells.Item(12, 8).Value = 2.250359
$ws.Cells.Item(12, 9).Value = 0.01049809664416703
$ws.Cells.Item(12, 10).Value = 0.01049809664416703
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 6.978882
$ws.Cells.Item(12, 14).Value = 20.936646
$ws.Cells.Item(12, 15).Value = 0.02499977909741928
$ws.Cells.Item(12, 16).Value = 0.02499977909741927
$ws.Cells.Item(12, 17).Value = 5.234996639545999
$ws.Cells.Item(12, 18).Value = 47.114969755914
$ws.Cells.Item(12, 19).Value = 0.0002624500970475345
$ws.Cells.Item(12, 20).Value = 0.0002624500970475345

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 0.7501196666666666
$ws.Cells.Item(13, 8).Value = 2.250359
$ws.Cells.Item(13, 9).Value = 0.01049809664416703
$ws.Cells.Item(13, 10).Value = 0.01049809664416703
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 7.105616
$ws.Cells.Item(13, 14).Value = 21.316848
$ws.Cells.Item(13, 15).Value = 0.02545376614063513
$ws.Cells.Item(13, 16).Value = 0.02545376614063513
$ws.Cells.Item(13, 17).Value = 5.330062305381333
$ws.Cells.Item(13, 18).Value = 47.970560748432
$ws.Cells.Item(13, 19).Value = 0.0002672160969024141
$ws.Cells.Item(13, 20).Value = 0.0002672160969024141

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 11.90182033333333
$ws.Cells.Item(14, 8).Value = 35.705461
$ws.Cells.Item(14, 9).Value = 0.166568703172488
$ws.Cells.Item(14, 10).Value = 0.166568703172488
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 218.7785543333333
$ws.Cells.Item(14, 14).Value = 656.3356630000001
$ws.Cells.Item(14, 15).Value = 0.7837094150017259
$ws.Cells.Item(14, 16).Value = 0.7837094150017259
$ws.Cells.Item(14, 17).Value = 2603.863046461738
$ws.Cells.Item(14, 18).Value = 23434.76741815565
$ws.Cells.Item(14, 19).Value = 0.1305414609209067
$ws.Cells.Item(14, 20).Value = 0.1305414609209067

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 11.90182033333333
$ws.Cells.Item(15, 8).Value = 35.705461
$ws.Cells.Item(15, 9).Value = 0.166568703172488
$ws.Cells.Item(15, 10).Value = 0.166568703172488
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 46.29469433333333
$ws.Cells.Item(15, 14).Value = 138.884083
$ws.Cells.Item(15, 15).Value = 0.1658370397602197
$ws.Cells.Item(15, 16).Value = 0.1658370397602197
$ws.Cells.Item(15, 17).Value = 550.991134341918
$ws.Cells.Item(15, 18).Value = 4958.920209077263
$ws.Cells.Item(15, 19).Value = 0.02762326065082412
$ws.Cells.Item(15, 20).Value = 0.02762326065082412

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 11.90182033333333
$ws.Cells.Item(16, 8).Value = 35.705461
$ws.Cells.Item(16, 9).Value = 0.166568703172488
$ws.Cells.Item(16, 10).Value = 0.166568703172488
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 6.978882
$ws.Cells.Item(16, 14).Value = 20.936646
$ws.Cells.Item(16, 15).Value = 0.02499977909741928
$ws.Cells.Item(16, 16).Value = 0.02499977909741927
$ws.Cells.Item(16, 17).Value = 83.06139969153399
$ws.Cells.Item(16, 18).Value = 747.5525972238059
$ws.Cells.Item(16, 19).Value = 0.004164180783855801
$ws.Cells.Item(16, 20).Value = 0.004164180783855801

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 11.90182033333333
$ws.Cells.Item(17, 8).Value = 35.705461
$ws.Cells.Item(17, 9).Value = 0.166568703172488
$ws.Cells.Item(17, 10).Value = 0.166568703172488
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 7.105616
$ws.Cells.Item(17, 14).Value = 21.316848
$ws.Cells.Item(17, 15).Value = 0.02545376614063513
$ws.Cells.Item(17, 16).Value = 0.02545376614063513
$ws.Cells.Item(17, 17).Value = 84.56976498965867
$ws.Cells.Item(17, 18).Value = 761.127884906928
$ws.Cells.Item(17, 19).Value = 0.004239800816901378
$ws.Cells.Item(17, 20).Value = 0.004239800816901377
